$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay a text value even when it looks like a plain
    # number (e.g. "212.26"), mirroring the original inline-string content.
    # Cells whose text contains multiple dots (e.g. "29.701.10") or other
    # non-numeric characters already round-trip as text without this.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.701.10"
$ws.Range("E2").Value = "  +3.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.604.28"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "212.26"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.19%  "

# Row 8 - Solana
Set-TextValue "D8" "27.24"
$ws.Range("E8").Value = "  +9.81%  "

# Row 9 - OKB
Set-TextValue "D9" "43.49"
$ws.Range("E9").Value = "  -1.32%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +2.10%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.14%  "

# Row 12 - TRON
Set-TextValue "D12" "0.0907"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.834.74"
$ws.Range("E13").Value = "  +2.62%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.596.65"
$ws.Range("E14").Value = "  +2.50%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "29.675.78"
$ws.Range("E15").Value = "  +3.54%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +4.02%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +2.40%  "

# Row 18 - Litecoin
Set-TextValue "D18" "63.35"
$ws.Range("E18").Value = "  +2.77%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "241.49"
$ws.Range("E19").Value = "  +6.07%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.97%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.43%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.10%  "

# Row 23 - Uniswap
Set-TextValue "D23" "4.00"
$ws.Range("E23").Value = "  +1.84%  "

# Row 24 - Avalanche
Set-TextValue "D24" "9.23"
$ws.Range("E24").Value = "  +1.89%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +0.92%  "

# Row 26 - Monero
Set-TextValue "D26" "154.96"
$ws.Range("E26").Value = "  +2.02%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +3.90%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.59%  "

# Row 29 - Cosmos
Set-TextValue "D29" "6.42"
$ws.Range("E29").Value = "  +2.60%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  -0.04%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +3.68%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.92%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.13%  "

# Row 34 - was InternetComputer(DFINITY), now Maker
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.431.12"
$ws.Range("E34").Value = "  +2.14%  "

# Row 35 - was Maker, now InternetComputer(DFINITY)
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D35" "3.12"
$ws.Range("E35").Value = "  +4.23%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +0.33%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +4.65%  "

# Row 38 - MXToken
Set-TextValue "D38" "2.81"
$ws.Range("E38").Value = "  +3.48%  "

# Row 39 - HuobiToken
$ws.Range("E39").Value = "  +0.24%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +1.44%  "

# Row 41 - ImmutableX
Set-TextValue "D41" "0.538"
$ws.Range("E41").Value = "  +4.12%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +2.32%  "

# Row 43 - BitcoinSV
Set-TextValue "D43" "54.64"
$ws.Range("E43").Value = "  +28.86%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  +5.83%  "

# Row 45 - ARBITRUM
Set-TextValue "D45" "0.801"
$ws.Range("E45").Value = "  +4.36%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  -0.09%  "

# Row 47 - Aave
Set-TextValue "D47" "65.85"
$ws.Range("E47").Value = "  +3.07%  "

# Row 48 - WEMIXToken
Set-TextValue "D48" "0.944"
$ws.Range("E48").Value = "  +11.79%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  +1.24%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "1.745.09"
$ws.Range("E50").Value = "  +2.79%  "

# Row 51 - Quant
Set-TextValue "D51" "86.56"
